$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Время обработки" (processing time) values in column C ---
$ws.Range("C2").Value = 0.006987
$ws.Range("C3").Value = 0.00163
$ws.Range("C4").Value = 0.013816
$ws.Range("C5").Value = 0.007046
$ws.Range("C6").Value = 0
$ws.Range("C9").Value = 0.005265
$ws.Range("C10").Value = 0.001656
$ws.Range("C15").Value = 0.005255
$ws.Range("C16").Value = 0.001503
$ws.Range("C21").Value = 0.00548
$ws.Range("C22").Value = 0.001067
$ws.Range("C23").Value = 0.001426
$ws.Range("C26").Value = 0.004445
$ws.Range("C27").Value = 0.001422
$ws.Range("C31").Value = 0.005514
$ws.Range("C32").Value = 0.001355
$ws.Range("C35").Value = 0
$ws.Range("C36").Value = 0.005565
$ws.Range("C37").Value = 0.001472
$ws.Range("C41").Value = 0.005496
$ws.Range("C42").Value = 0.001442
$ws.Range("C44").Value = 0
$ws.Range("C45").Value = 0.005487
$ws.Range("C46").Value = 0.00107
$ws.Range("C47").Value = 0.001514
$ws.Range("C48").Value = 0.005931
$ws.Range("C49").Value = 0.006471
$ws.Range("C50").Value = 0.000578
$ws.Range("C51").Value = 0.00685
$ws.Range("C52").Value = 0.006948
$ws.Range("C53").Value = 0.006896
$ws.Range("C54").Value = 0.005375
$ws.Range("C55").Value = 0.001527
$ws.Range("C56").Value = 0.006982
$ws.Range("C57").Value = 0.006978
$ws.Range("C58").Value = 0.006931
$ws.Range("C59").Value = 0.005412
$ws.Range("C60").Value = 0.00045
$ws.Range("C61").Value = 0.006912
$ws.Range("C62").Value = 0.006934
$ws.Range("C63").Value = 0.006975
$ws.Range("C64").Value = 0.006465

# --- Append 20 new "black_fred_*" rows (65-84) ---
$ws.Range("B65").NumberFormat = "@"
$ws.Range("A65").Value = "black_fred_1.jpg"
$ws.Range("B65").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C65").Value = 0.006916
$ws.Range("D65").Value = 0

$ws.Range("B66").NumberFormat = "@"
$ws.Range("A66").Value = "black_fred_2.jpg"
$ws.Range("B66").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C66").Value = 0.006896
$ws.Range("D66").Value = 0

$ws.Range("B67").NumberFormat = "@"
$ws.Range("A67").Value = "black_fred_3.jpg"
$ws.Range("B67").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C67").Value = 0.006519
$ws.Range("D67").Value = 0

$ws.Range("B68").NumberFormat = "@"
$ws.Range("A68").Value = "black_fred_4.jpg"
$ws.Range("B68").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C68").Value = 0.007991
$ws.Range("D68").Value = 0

$ws.Range("B69").NumberFormat = "@"
$ws.Range("A69").Value = "black_fred_5.jpg"
$ws.Range("B69").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C69").Value = 0.005639
$ws.Range("D69").Value = 0

$ws.Range("B70").NumberFormat = "@"
$ws.Range("A70").Value = "black_fred_6.jpg"
$ws.Range("B70").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C70").Value = 0.00572
$ws.Range("D70").Value = 0

$ws.Range("B71").NumberFormat = "@"
$ws.Range("A71").Value = "black_fred_7.jpg"
$ws.Range("B71").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C71").Value = 0.007012
$ws.Range("D71").Value = 0

$ws.Range("B72").NumberFormat = "@"
$ws.Range("A72").Value = "black_fred_8.jpg"
$ws.Range("B72").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C72").Value = 0.005847
$ws.Range("D72").Value = 0

$ws.Range("B73").NumberFormat = "@"
$ws.Range("A73").Value = "black_fred_9.jpg"
$ws.Range("B73").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C73").Value = 0.00592
$ws.Range("D73").Value = 0

$ws.Range("B74").NumberFormat = "@"
$ws.Range("A74").Value = "black_fred_10.jpg"
$ws.Range("B74").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C74").Value = 0.001662
$ws.Range("D74").Value = 0

$ws.Range("B75").NumberFormat = "@"
$ws.Range("A75").Value = "black_fred_11.jpg"
$ws.Range("B75").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C75").Value = 0.006885
$ws.Range("D75").Value = 0

$ws.Range("B76").NumberFormat = "@"
$ws.Range("A76").Value = "black_fred_12.jpg"
$ws.Range("B76").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C76").Value = 0.007135
$ws.Range("D76").Value = 0

$ws.Range("B77").NumberFormat = "@"
$ws.Range("A77").Value = "black_fred_13.jpg"
$ws.Range("B77").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C77").Value = 0.006697
$ws.Range("D77").Value = 0

$ws.Range("B78").NumberFormat = "@"
$ws.Range("A78").Value = "black_fred_14.jpg"
$ws.Range("B78").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C78").Value = 0.005392
$ws.Range("D78").Value = 0

$ws.Range("B79").NumberFormat = "@"
$ws.Range("A79").Value = "black_fred_15.jpg"
$ws.Range("B79").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C79").Value = 0.002655
$ws.Range("D79").Value = 0

$ws.Range("B80").NumberFormat = "@"
$ws.Range("A80").Value = "black_fred_16.jpg"
$ws.Range("B80").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C80").Value = 0.005933
$ws.Range("D80").Value = 0

$ws.Range("B81").NumberFormat = "@"
$ws.Range("A81").Value = "black_fred_17.jpg"
$ws.Range("B81").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C81").Value = 0.007236
$ws.Range("D81").Value = 0

$ws.Range("B82").NumberFormat = "@"
$ws.Range("A82").Value = "black_fred_18.jpg"
$ws.Range("B82").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C82").Value = 0.006802
$ws.Range("D82").Value = 0

$ws.Range("B83").NumberFormat = "@"
$ws.Range("A83").Value = "black_fred_19.jpg"
$ws.Range("B83").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C83").Value = 0.006187
$ws.Range("D83").Value = 0

$ws.Range("B84").NumberFormat = "@"
$ws.Range("A84").Value = "black_fred_20.jpg"
$ws.Range("B84").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C84").Value = 0.005924
$ws.Range("D84").Value = 0

